$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. The "Last updated" date placeholder (a datetimeFigureOut field) on the
#    slide master and every slide layout is bumped from 9/5/2013 to
#    4/28/2014.
# ---------------------------------------------------------------------------
$m = $p.SlideMaster

for ($k = 1; $k -le $m.Shapes.Count; $k++) {
    $msh = $m.Shapes.Item($k)
    if ($msh.Name -like "Date Placeholder*") {
        $msh.TextFrame.TextRange.Text = "4/28/2014"
    }
}

$cl = $m.CustomLayouts
for ($i = 1; $i -le $cl.Count; $i++) {
    $layout = $cl.Item($i)
    $lshapes = $layout.Shapes
    for ($j = 1; $j -le $lshapes.Count; $j++) {
        $lsh = $lshapes.Item($j)
        if ($lsh.Name -like "Date Placeholder*") {
            $lsh.TextFrame.TextRange.Text = "4/28/2014"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 1 gets a second copy of the "Bush-dog.jpg" picture (the start of
#    captioned-image support), placed to the right of the title area.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$srcPic = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $cand = $s1.Shapes.Item($i)
    if ($cand.AlternativeText -like "*Bush-dog.jpg") {
        $srcPic = $cand
    }
}

# Burn shape-id slots 2 and 3 with cheap throw-away textboxes so that the
# picture duplicated below lands on id 6 -- the id it was given originally.
$burn1 = $s1.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn1.Delete()
$burn2 = $s1.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn2.Delete()

$newPicRange = $srcPic.Duplicate()
$newPic = $newPicRange.Item(1)
$newPic.Left = 492
$newPic.Top = 57.0536
$newPic.Width = 144
$newPic.Height = 143.25
